# Edit the workbook per commit "#5: fund, bonds, otherbonds, antique done"
#
# 1. Remove the "其他有價證券" (other securities) sheet entirely.
# 2. Rebuild the "基金受益憑證" (fund) sheet with the full record schema
#    (name/owner/dealer/quantity/face_value/currency/total/property_category/
#    category/date/legislator_name/legislator_id/source_file/index), matching
#    the layout already used by the "股票" (stock) sheet, plus a "dealer"
#    column for the fund distributor/bank.

$wb = $excel.ActiveWorkbook

# --- 1. Drop the "其他有價證券" sheet -------------------------------------
$wb.Worksheets.Item("其他有價證券").Delete()

# --- 2. Rebuild "基金受益憑證" -------------------------------------------
$ws = $wb.Worksheets.Item("基金受益憑證")

# Extend formatting (bold/border header style, plain data style) from the
# existing B:H columns into the new I:O columns before filling values in,
# so the new cells pick up the same cell styles as their row neighbours.
$ws.Range("B1:H1").Copy()
$ws.Range("I1:O1").PasteSpecial(-4122)

$ws.Range("B2:H4").Copy()
$ws.Range("I2:O4").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Header row
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "dealer"
$ws.Range("E1").Value = "quantity"
$ws.Range("F1").Value = "face_value"
$ws.Range("G1").Value = "currency"
$ws.Range("H1").Value = "total"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

# Row 2 - 德利全球資源產業基金
$ws.Range("A2").Value = 71
$ws.Range("B2").Value = "德利全球資源產業基金"
$ws.Range("C2").Value = "田秋堇"
$ws.Range("D2").Value = "台北富邦商業銀行"
$ws.Range("E2").Value = 69.772
$ws.Range("F2").Value = 3835
$ws.Range("G2").Value = "新臺幣"
$ws.Range("H2").Value = 267575.62
$ws.Range("I2").Value = "fund"
$ws.Range("J2").Value = "normal"
$ws.Range("K2").Value = "2012-04-10"
$ws.Range("L2").Value = "田秋堇"
$ws.Range("M2").Value = 1316
$ws.Range("N2").Value = "tmp9b251"
$ws.Range("O2").Value = 71

# Row 3 - 聯博美國收益澳幣避險基金
$ws.Range("A3").Value = 72
$ws.Range("B3").Value = "聯博美國收益澳幣避險基金"
$ws.Range("C3").Value = "田秋堇"
$ws.Range("D3").Value = "台北富邦商業銀行"
$ws.Range("E3").Value = 5315.604
$ws.Range("F3").Value = 464
$ws.Range("G3").Value = "新臺幣"
$ws.Range("H3").Value = 2466440.26
$ws.Range("I3").Value = "fund"
$ws.Range("J3").Value = "normal"
$ws.Range("K3").Value = "2012-04-10"
$ws.Range("L3").Value = "田秋堇"
$ws.Range("M3").Value = 1316
$ws.Range("N3").Value = "tmp9b251"
$ws.Range("O3").Value = 72

# Row 4 - 聯博美國收益澳幣避險基金 (second lot, AUD hedge share)
$ws.Range("A4").Value = 73
$ws.Range("B4").Value = "聯博美國收益澳幣避險基金"
$ws.Range("C4").Value = "田秋堇"
$ws.Range("D4").Value = "台北富邦商業銀行"
$ws.Range("E4").Value = 648.508
$ws.Range("F4").Value = 15.33
$ws.Range("G4").Value = "澳幣"
$ws.Range("H4").Value = 301931.21
$ws.Range("I4").Value = "fund"
$ws.Range("J4").Value = "normal"
$ws.Range("K4").Value = "2012-04-10"
$ws.Range("L4").Value = "田秋堇"
$ws.Range("M4").Value = 1316
$ws.Range("N4").Value = "tmp9b251"
$ws.Range("O4").Value = 73
